$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.250.04"
$ws.Range("E2").Value = "  +6.64%  "
$ws.Range("D3").Value = "3.017.52"
$ws.Range("E3").Value = "  +3.50%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'585.47"
$ws.Range("D6").Value = "'163.30"
$ws.Range("E6").Value = "  +13.22%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.013.10"
$ws.Range("E8").Value = "  +3.38%  "
$ws.Range("E9").Value = "  +3.83%  "
$ws.Range("D10").Value = "'6.73"
$ws.Range("E10").Value = "  -3.52%  "
$ws.Range("D11").Value = "'0.157"
$ws.Range("E11").Value = "  +5.48%  "
$ws.Range("E12").Value = "  +6.22%  "
$ws.Range("D13").Value = "'0.0000258"
$ws.Range("E13").Value = "  +8.22%  "
$ws.Range("D14").Value = "'34.81"
$ws.Range("E14").Value = "  +6.53%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "66.167.02"
$ws.Range("E16").Value = "  +6.58%  "
$ws.Range("D17").Value = "3.519.72"
$ws.Range("E18").Value = "  +6.77%  "
$ws.Range("D19").Value = "3.019.29"
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("D20").Value = "'458.98"
$ws.Range("E20").Value = "  +6.26%  "
$ws.Range("D21").Value = "'13.94"
$ws.Range("E21").Value = "  +6.52%  "
$ws.Range("D22").Value = "'0.690"
$ws.Range("E22").Value = "  +5.56%  "
$ws.Range("D23").Value = "'7.41"
$ws.Range("E23").Value = "  +7.91%  "
$ws.Range("D24").Value = "'82.46"
$ws.Range("E24").Value = "  +4.43%  "
$ws.Range("D25").Value = "'2.31"
$ws.Range("E25").Value = "  +14.39%  "
$ws.Range("D26").Value = "'12.41"
$ws.Range("E26").Value = "  +3.08%  "
$ws.Range("D27").Value = "'10.57"
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "'8.09"
$ws.Range("E29").Value = "  +16.17%  "
$ws.Range("D30").Value = "'2.38"
$ws.Range("E30").Value = "  +16.99%  "
$ws.Range("E31").Value = "  -7.11%  "
$ws.Range("E32").Value = "  +4.38%  "
$ws.Range("D33").Value = "'27.37"
$ws.Range("E33").Value = "  +6.36%  "
$ws.Range("E34").Value = "  +4.97%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'0.994"
$ws.Range("E36").Value = "  +3.83%  "
$ws.Range("D37").Value = "'5.83"
$ws.Range("E37").Value = "  +7.88%  "
$ws.Range("E38").Value = "  +15.83%  "
$ws.Range("E39").Value = "  +3.87%  "
$ws.Range("D40").Value = "'50.06"
$ws.Range("D41").Value = "'0.309"
$ws.Range("E41").Value = "  +15.73%  "
$ws.Range("E42").Value = "  +8.14%  "
$ws.Range("D43").Value = "'43.67"
$ws.Range("E43").Value = "  +5.93%  "
$ws.Range("E44").Value = "  +3.67%  "
$ws.Range("D45").Value = "'397.05"
$ws.Range("E45").Value = "  +13.79%  "
$ws.Range("D46").Value = "'0.0363"
$ws.Range("E46").Value = "  +7.53%  "
$ws.Range("D47").Value = "2.805.11"
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").Value = "'134.21"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'24.03"
$ws.Range("E50").Value = "  +11.17%  "
$ws.Range("E51").Value = "  +4.32%  "
